$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new test case row (row 10): "Update user profile with summary"
$ws.Range("A10").Value = "S1_TC_T9"
$ws.Range("B10").Value = "Update user profile with summary"
$ws.Range("C10").Value = "1PPROFILE"
$ws.Range("D10").Value = "/users/user/(SYS_USER1)"
$ws.Range("E10").Value = "PUT"
$ws.Range("F10").Value = "Content-Type=application/json"
$ws.Range("H10").Value = '{"summary":"test"}'
$ws.Range("J10").Value = "status=200||summary=test"
$ws.Range("L10").Value = "PASS"

# Materialize the empty cells in columns G, I, K for row 10 (to mirror other rows)
$ws.Range("G10").Borders.LineStyle = -4142
$ws.Range("I10").Borders.LineStyle = -4142
$ws.Range("K10").Borders.LineStyle = -4142

# J column has no column-level wrap style, so apply it explicitly to J10
$ws.Range("J10").WrapText = $true

# Column H width adjustment (widened to fit the new body text)
$ws.Columns.Item(8).ColumnWidth = 79.67

# Update the view state to match the new selection/scroll position
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("H10").Select()
